# SSDM-12286 Fixed letter case inconsistencies.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Vocabulary Code" -> "Vocabulary code"
$ws.Range("H4").Value = "Vocabulary code"

# "Generated Code Prefix" -> "Generated code prefix"
$ws.Range("E2").Value = "Generated code prefix"

# Update the selection to match target workbook state
$ws.Range("E3").Select()
